# Updates Price (D) and Volume(1h) (E) columns on the cryptos sheet
# with refreshed values, matching the "Updated cryptos list" GitHub Actions commit.
# Price values that look like plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (preserving exact formatting,
# e.g. trailing zeros / no float rounding) instead of coercing to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.724.35"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.646.65"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'213.26"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D8").Value = "'23.34"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.879.03"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.634.27"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "'0.560"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'64.76"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "27.700.43"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").Value = "'231.91"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "'7.64"
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'4.30"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'10.17"
$ws.Range("E23").Value = "  +9.48%  "
$ws.Range("E24").Value = "  -3.53%  "
$ws.Range("D25").Value = "'150.17"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "'6.92"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("D28").Value = "'15.66"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").Value = "'0.0487"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "1.441.11"
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("D37").Value = "'0.571"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").Value = "'0.879"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "'0.887"
$ws.Range("E40").Value = "  +12.64%  "
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "'67.20"
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("D44").Value = "'5.58"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").Value = "1.788.06"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "'1.75"
$ws.Range("E47").Value = "  +6.43%  "
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").Value = "'85.60"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "'0.0989"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "'7.77"
$ws.Range("E51").Value = "  +2.15%  "
